# Update the "想去人数" (want-to-go count) figures in the F column
# on the "展览" and "全部类型" worksheets, reflecting new data pulled
# at build time (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1050
$wsExpo.Range("F5").Value = 2878
$wsExpo.Range("F11").Value = 126
$wsExpo.Range("F12").Value = 46
$wsExpo.Range("F13").Value = 2706
$wsExpo.Range("F14").Value = 955

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1050
$wsAll.Range("F6").Value = 2878
$wsAll.Range("F13").Value = 126
$wsAll.Range("F14").Value = 46
$wsAll.Range("F15").Value = 2706
$wsAll.Range("F16").Value = 955
